$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 ("ajouter regle: Brelan") test results:
# D9 "Qui" -> Liu, Huang
$ws.Range("D9").Value = "Liu, Huang"

# I9 "Fin" -> 2018-03-08 (date serial 43167), matching the date format
# already used by H9 ("Debut") / I8. Copy the number format from H9
# first so we reuse the existing date style instead of creating a new one.
$ws.Range("H9").Copy()
$ws.Range("I9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("I9").Value = 43167

# J9 "Test ?" -> OK
$ws.Range("J9").Value = "OK"

# Reflect the author's last selection in the saved sheet view
$ws.Range("E11").Select() | Out-Null
